$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Sucrose__Phellogen_Day_sp_exchange", -0.02275250248137701, -0.02573249889733798, 0.02524392192411239, $false),
    @("Nitrate__Ibark_Day_sp_exchange", -0.2094768778078743, -0.2095535822389241, -0.2094694318056584, $true),
    @("L-Leucine__Leaf_Day_sp_exchange", -0.002733808561680673, -0.00288651637994059, -0.002709062279195087, $true),
    @("L-Leucine__Ibark_Day_sp_exchange", -0.005245730614982, -0.006362180098646557, -0.004125825124105065, $true),
    @("L-Threonine__Leaf_Day_sp_exchange", -0.006957364279378033, -0.007034642068470099, -0.006944911409846659, $true),
    @("L-Valine__Leaf_Day_sp_exchange", -0.008344632411180001, -0.008362952132431756, -0.008337433269960105, $true),
    @("L-Valine__Ibark_Day_sp_exchange", -0.003496901480776, -0.004241147259572694, -0.002749189943002829, $true),
    @("L-Glutamine__Leaf_Day_sp_exchange", -0.005197582981536614, -0.005197582979205357, -0.005197582985947076, $true),
    @("L-Glutamine__Ibark_Day_sp_exchange", -0.006012063502197787, -0.01023537158539439, -0.001788755371801476, $true),
    @("Glycine__Ibark_Day_sp_exchange", -0.1800205736429344, -0.1812062262971476, -0.007306397149129167, $true),
    @("Glycine__Phellogen_Day_sp_exchange", -0.01642057489156097, -0.01711900264407966, -0.01640512034329712, $true),
    @("L-Glutamate__Leaf_Day_sp_exchange", -0.3281353227496608, -0.3281860887639729, -0.3281038887883893, $true),
    @("L-Glutamate__Phellogen_Day_sp_exchange", -0.006748869194980841, -0.1071260099573368, 0.002446565458976557, $false),
    @("L-Asparagine__Leaf_Day_sp_exchange", -0.005459338897155001, -0.005461144685913809, -0.005454830852780226, $true),
    @("L-Asparagine__Ibark_Day_sp_exchange", -0.002287790442146, -0.00277470053989251, -0.001800880345507929, $true),
    @("L-Serine__Leaf_Day_sp_exchange", -0.09109761126615967, -0.250490096289229, -0.0397291287589565, $true),
    @("Citrate__Leaf_Day_sp_exchange", -0.006637622814004108, -0.006676289670016933, -0.006596423716057398, $true)
)

$formatSource = $ws.Cells.Item(2, 1)
$formatSource.Copy()

$rowIndex = 2
foreach ($item in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $item[0]
    $ws.Cells.Item($rowIndex, 2).Value = $item[1]
    $ws.Cells.Item($rowIndex, 3).Value = $item[2]
    $ws.Cells.Item($rowIndex, 4).Value = $item[3]
    $ws.Cells.Item($rowIndex, 5).Value = $item[4]
    if ($rowIndex -gt 2) {
        $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    }
    $rowIndex++
}

$excel.CutCopyMode = $false
